$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new time entry for row 14
$ws.Range("C14").Value = 1.5
$ws.Range("D14").Value = "New branch for old Helsinkikanava API."

# Update the view: scroll/selection
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D15").Select()
